$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 2019-07-07 (row 8): H8 flips from "√" to "×"
$ws.Range("H8").Value = "×"

# 2019-07-08 (row 9): fill in the full day's checklist
$ws.Range("B9").Value = "√"
$ws.Range("C9").Value = "×"
$ws.Range("D9").Value = "√"
$ws.Range("E9").Value = "√"
$ws.Range("F9").Value = "√"
$ws.Range("G9").Value = "√"
$ws.Range("H9").Value = "×"
$ws.Range("I9").Value = "√"
$ws.Range("J9").Value = "√"
$ws.Range("K9").Value = "√"
$ws.Range("L9").Value = "√"
$ws.Range("M9").Value = "√"
$ws.Range("N9").Value = "√"
$ws.Range("O9").Value = "√"

# 2019-07-09 (row 10): new day started, only the date so far
$ws.Range("A10").Value = 20190709

# Mirror the cursor move the author made after typing the new date
$ws.Range("B10").Select()
